$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4, pushing existing rows 4-10 down to 5-11.
$ws.Rows("4:4").Insert()

# Fill the new row 4 with the new record (same pattern as surrounding rows,
# new date and price figures).
$ws.Range("A4").Value = 4
$ws.Range("B4").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C4").Value = "Los Lagos"
$ws.Range("D4").Value = 44764
$ws.Range("E4").Value = 10
$ws.Range("F4").Value = 100112035
$ws.Range("G4").Value = "Bruselas (repollito)"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 90
$ws.Range("K4").Value = 24000
$ws.Range("L4").Value = 24000
$ws.Range("M4").Value = 24000
$ws.Range("N4").Value = "$/malla 15 kilos"
$ws.Range("O4").Value = "Provincia de Quillota"
$ws.Range("P4").Value = 1600
$ws.Range("Q4").Value = 15
$ws.Range("R4").Value = "Hortaliza"
